$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every D/E cell in this sheet is stored as plain text (t="inlineStr" in the
# original OOXML) even though most of the values look numeric. Writing a
# numeric-looking string straight into Range.Value lets Excel auto-convert it
# to a real number (dropping significant trailing zeros, e.g. "0.3860" ->
# 0.386, and flipping the stored cell type). Prefixing the literal with a
# leading apostrophe forces Excel to keep it as text; re-applying the
# "Normal" style right after removes the quote-prefix cell format Excel
# stamps on such cells, so no visible/formatting side effect remains.
$textCells = New-Object System.Collections.ArrayList

$ws.Range("D2").Value = "'28.009.13"
[void]$textCells.Add("D2")
$ws.Range("E2").Value = "'  -2.31%  "
[void]$textCells.Add("E2")
$ws.Range("D3").Value = "'1.829.11"
[void]$textCells.Add("D3")
$ws.Range("E3").Value = "'  -1.40%  "
[void]$textCells.Add("E3")
$ws.Range("D4").Value = "'1.002"
[void]$textCells.Add("D4")
$ws.Range("E4").Value = "'  -0.29%  "
[void]$textCells.Add("E4")
$ws.Range("D5").Value = "'326.69"
[void]$textCells.Add("D5")
$ws.Range("E5").Value = "'  -2.82%  "
[void]$textCells.Add("E5")
$ws.Range("D6").Value = "'1.002"
[void]$textCells.Add("D6")
$ws.Range("E6").Value = "'  -0.22%  "
[void]$textCells.Add("E6")
$ws.Range("D7").Value = "'0.4614"
[void]$textCells.Add("D7")
$ws.Range("E7").Value = "'  -0.83%  "
[void]$textCells.Add("E7")
$ws.Range("D8").Value = "'0.3860"
[void]$textCells.Add("D8")
$ws.Range("E8").Value = "'  -1.64%  "
[void]$textCells.Add("E8")
$ws.Range("D9").Value = "'0.07859"
[void]$textCells.Add("D9")
$ws.Range("E9").Value = "'  -0.92%  "
[void]$textCells.Add("E9")
$ws.Range("D10").Value = "'0.9574"
[void]$textCells.Add("D10")
$ws.Range("E10").Value = "'  -2.75%  "
[void]$textCells.Add("E10")
$ws.Range("D11").Value = "'21.83"
[void]$textCells.Add("D11")
$ws.Range("E11").Value = "'  -1.30%  "
[void]$textCells.Add("E11")
$ws.Range("D12").Value = "'1.855.20"
[void]$textCells.Add("D12")
$ws.Range("E12").Value = "'  -0.46%  "
[void]$textCells.Add("E12")
$ws.Range("D13").Value = "'5.650"
[void]$textCells.Add("D13")
$ws.Range("E13").Value = "'  -3.42%  "
[void]$textCells.Add("E13")
$ws.Range("D14").Value = "'6.875"
[void]$textCells.Add("D14")
$ws.Range("E14").Value = "'  -2.16%  "
[void]$textCells.Add("E14")
$ws.Range("D15").Value = "'0.06750"
[void]$textCells.Add("D15")
$ws.Range("E15").Value = "'  -0.69%  "
[void]$textCells.Add("E15")
$ws.Range("D16").Value = "'1.002"
[void]$textCells.Add("D16")
$ws.Range("E16").Value = "'  -0.38%  "
[void]$textCells.Add("E16")
$ws.Range("D17").Value = "'86.74"
[void]$textCells.Add("D17")
$ws.Range("E17").Value = "'  -1.14%  "
[void]$textCells.Add("E17")
$ws.Range("D18").Value = "'0.000009915"
[void]$textCells.Add("D18")
$ws.Range("E18").Value = "'  -2.26%  "
[void]$textCells.Add("E18")
$ws.Range("D19").Value = "'16.59"
[void]$textCells.Add("D19")
$ws.Range("E19").Value = "'  -2.76%  "
[void]$textCells.Add("E19")
$ws.Range("E20").Value = "'  -0.34%  "
[void]$textCells.Add("E20")
$ws.Range("D21").Value = "'28.037.46"
[void]$textCells.Add("D21")
$ws.Range("D22").Value = "'5.301"
[void]$textCells.Add("D22")
$ws.Range("E22").Value = "'  -2.16%  "
[void]$textCells.Add("E22")
$ws.Range("D23").Value = "'10.97"
[void]$textCells.Add("D23")
$ws.Range("E23").Value = "'  -3.31%  "
[void]$textCells.Add("E23")
$ws.Range("D24").Value = "'2.089"
[void]$textCells.Add("D24")
$ws.Range("E24").Value = "'  -2.11%  "
[void]$textCells.Add("E24")
$ws.Range("D25").Value = "'2.111.63"
[void]$textCells.Add("D25")
$ws.Range("E25").Value = "'  +0.83%  "
[void]$textCells.Add("E25")
$ws.Range("D26").Value = "'153.71"
[void]$textCells.Add("D26")
$ws.Range("E26").Value = "'  +0.26%  "
[void]$textCells.Add("E26")
$ws.Range("D27").Value = "'19.16"
[void]$textCells.Add("D27")
$ws.Range("E27").Value = "'  -1.67%  "
[void]$textCells.Add("E27")
$ws.Range("D28").Value = "'5.719"
[void]$textCells.Add("D28")
$ws.Range("E28").Value = "'  -8.67%  "
[void]$textCells.Add("E28")
$ws.Range("D29").Value = "'1.968"
[void]$textCells.Add("D29")
$ws.Range("E29").Value = "'  -3.20%  "
[void]$textCells.Add("E29")
$ws.Range("D30").Value = "'117.15"
[void]$textCells.Add("D30")
$ws.Range("E30").Value = "'  -0.44%  "
[void]$textCells.Add("E30")
$ws.Range("D31").Value = "'0.9346"
[void]$textCells.Add("D31")
$ws.Range("E31").Value = "'  -4.97%  "
[void]$textCells.Add("E31")
$ws.Range("D32").Value = "'0.09237"
[void]$textCells.Add("D32")
$ws.Range("E32").Value = "'  -2.26%  "
[void]$textCells.Add("E32")
$ws.Range("D33").Value = "'5.287"
[void]$textCells.Add("D33")
$ws.Range("E33").Value = "'  -2.12%  "
[void]$textCells.Add("E33")
$ws.Range("D34").Value = "'1.313"
[void]$textCells.Add("D34")
$ws.Range("E34").Value = "'  -2.62%  "
[void]$textCells.Add("E34")
$ws.Range("D35").Value = "'3.323"
[void]$textCells.Add("D35")
$ws.Range("E35").Value = "'  -5.31%  "
[void]$textCells.Add("E35")
$ws.Range("D36").Value = "'0.05855"
[void]$textCells.Add("D36")
$ws.Range("E36").Value = "'  -4.94%  "
[void]$textCells.Add("E36")
$ws.Range("D37").Value = "'0.02140"
[void]$textCells.Add("D37")
$ws.Range("E37").Value = "'  -2.76%  "
[void]$textCells.Add("E37")
$ws.Range("E38").Value = "'  -1.36%  "
[void]$textCells.Add("E38")
$ws.Range("E39").Value = "'  +1.40%  "
[void]$textCells.Add("E39")
$ws.Range("D40").Value = "'0.5575"
[void]$textCells.Add("D40")
$ws.Range("E40").Value = "'  -2.80%  "
[void]$textCells.Add("E40")
$ws.Range("D41").Value = "'9.866"
[void]$textCells.Add("D41")
$ws.Range("E41").Value = "'  -2.45%  "
[void]$textCells.Add("E41")
$ws.Range("D42").Value = "'0.1757"
[void]$textCells.Add("D42")
$ws.Range("E42").Value = "'  -1.80%  "
[void]$textCells.Add("E42")
$ws.Range("D43").Value = "'1.225"
[void]$textCells.Add("D43")
$ws.Range("E43").Value = "'  -2.19%  "
[void]$textCells.Add("E43")
$ws.Range("D44").Value = "'11.56"
[void]$textCells.Add("D44")
$ws.Range("E44").Value = "'  -2.85%  "
[void]$textCells.Add("E44")
$ws.Range("D45").Value = "'0.5261"
[void]$textCells.Add("D45")
$ws.Range("E45").Value = "'  -2.83%  "
[void]$textCells.Add("E45")
$ws.Range("D46").Value = "'0.07027"
[void]$textCells.Add("D46")
$ws.Range("E46").Value = "'  -1.56%  "
[void]$textCells.Add("E46")
$ws.Range("D47").Value = "'2.134"
[void]$textCells.Add("D47")
$ws.Range("E47").Value = "'  -9.33%  "
[void]$textCells.Add("E47")
$ws.Range("D48").Value = "'1.824"
[void]$textCells.Add("D48")
$ws.Range("E48").Value = "'  -4.67%  "
[void]$textCells.Add("E48")
$ws.Range("D49").Value = "'112.49"
[void]$textCells.Add("D49")
$ws.Range("E49").Value = "'  -2.73%  "
[void]$textCells.Add("E49")
$ws.Range("D50").Value = "'1.001"
[void]$textCells.Add("D50")
$ws.Range("E50").Value = "'  -0.42%  "
[void]$textCells.Add("E50")
$ws.Range("D51").Value = "'2.320"
[void]$textCells.Add("D51")
$ws.Range("E51").Value = "'  -0.28%  "
[void]$textCells.Add("E51")

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

